$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $text)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '69.673.56'
$ws.Range("E2").Value = '  +2.80%  '
Set-TextValue "D3" '2.516.04'
$ws.Range("E3").Value = '  +0.92%  '
$ws.Range("E4").Value = '  +0.04%  '
Set-TextValue "D5" '598.24'
$ws.Range("E5").Value = '  +1.89%  '
Set-TextValue "D6" '176.64'
$ws.Range("E6").Value = '  -0.27%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E8").Value = '  +0.88%  '
Set-TextValue "D9" '2.515.17'
$ws.Range("E9").Value = '  +0.93%  '
Set-TextValue "D10" '0.158'
$ws.Range("E10").Value = '  +12.37%  '
$ws.Range("E11").Value = '  -0.32%  '
Set-TextValue "D12" '0.344'
$ws.Range("E12").Value = '  +1.19%  '
$ws.Range("E13").Value = '  +1.44%  '
Set-TextValue "D14" '2.977.34'
$ws.Range("E14").Value = '  +1.02%  '
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue "D15" '0.0000180'
$ws.Range("E15").Value = '  +5.15%  '
$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextValue "D16" '25.94'
$ws.Range("E16").Value = '  +1.16%  '
Set-TextValue "D17" '69.533.76'
$ws.Range("E17").Value = '  +2.84%  '
Set-TextValue "D18" '2.496.83'
$ws.Range("E18").Value = '  -0.16%  '
Set-TextValue "D19" '7.65'
Set-TextValue "D20" '363.05'
$ws.Range("E20").Value = '  +3.14%  '
Set-TextValue "D21" '11.01'
$ws.Range("E21").Value = '  +0.52%  '
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("E23").Value = '  -0.09%  '
Set-TextValue "D24" '70.60'
$ws.Range("E24").Value = '  -0.29%  '
$ws.Range("E25").Value = '  -1.18%  '
$ws.Range("E26").Value = '  -1.23%  '
Set-TextValue "D27" '9.12'
$ws.Range("E27").Value = '  +0.16%  '
Set-TextValue "D28" '2.644.50'
$ws.Range("E28").Value = '  +2.18%  '
Set-TextValue "D29" '0.999'
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D30" '0.0₃0897'
$ws.Range("E30").Value = '  -1.05%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue "D31" '509.67'
$ws.Range("E31").Value = '  +1.00%  '
Set-TextValue "D32" '7.71'
$ws.Range("E32").Value = '  -1.71%  '
$ws.Range("E33").Value = '  -1.47%  '
$ws.Range("E34").Value = '  +0.80%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  -1.35%  '
Set-TextValue "D37" '160.19'
$ws.Range("E37").Value = '  -2.34%  '
Set-TextValue "D38" '18.73'
$ws.Range("E38").Value = '  +2.14%  '
Set-TextValue "D39" '18.91'
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("E40").Value = '  -1.18%  '
$ws.Range("E41").Value = '  +0.05%  '
Set-TextValue "D42" '1.72'
$ws.Range("E42").Value = '  -0.94%  '
$ws.Range("E43").Value = '  -1.38%  '
Set-TextValue "D44" '0.322'
$ws.Range("E44").Value = '  -2.31%  '
Set-TextValue "D45" '2.37'
$ws.Range("E45").Value = '  -3.33%  '
Set-TextValue "D46" '38.80'
$ws.Range("E46").Value = '  -0.39%  '
Set-TextValue "D47" '150.55'
$ws.Range("E47").Value = '  +4.12%  '
$ws.Range("E48").Value = '  +1.76%  '
$ws.Range("E49").Value = '  +0.13%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D50" '0.0738'
$ws.Range("E50").Value = '  -0.67%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue "D51" '0.0₆0251'
$ws.Range("E51").Value = '  -1.42%  '

Write-Host "Applied cryptos update"
